$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.261.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.447.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.73%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9421'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.99%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '274.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3064'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.56'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.52%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.027'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06505'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.35%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9968'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.343'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.065'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001011'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.442.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9557'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05685'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.388'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.58%  '

$ws.Range("E23").Value = '  -2.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.237'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.275.28'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.075'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.593.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '110.81'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.948'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.74%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7854'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.26%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.784'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.09%  '

$ws.Range("E34").Value = '  +0.20%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.466'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05651'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.93%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.645'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.115'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.34%  '

$ws.Range("B39").Value = 'Frax'
$ws.Range("C39").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9467'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.59%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02000'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1846'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.154'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -15.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5218'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.40%  '

$ws.Range("E45").Value = '  -1.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.83'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '116.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5096'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.731'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06372'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9799'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.36%  '
